$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New SMA connector (part added to BOM, but not yet to the schematic):
# replace the old part number / cost on the "SMA connector" BOM row (row 15)
# with the new vendor part number and its unit cost.
$ws.Range("B15").Value = "A97594-ND"
$ws.Range("E15").Value = 2.17

# Mirror the author's final selection/view state after the edit.
$ws.Range("E16").Select()
